$wb = $excel.ActiveWorkbook

# --- Zub_Gel: insert a "Rabatt in %" (p_rabatt) parameter row and update the
#     final price formula to apply the discount ----------------------------
$wsZubGel = $wb.Sheets.Item("Zub_Gel")

# Push the old "Endpreis" row (row 8) down to make room for the new row.
[void]$wsZubGel.Rows.Item(8).Insert()

$wsZubGel.Range("A8").Value = "Zahl"
$wsZubGel.Range("B8").Value = "Rabatt in %"
$wsZubGel.Range("C8").Value = "p_rabatt"
$wsZubGel.Range("D8").Value = 15

# Update the formula description cell (now on row 9) so it applies the
# discount, matching the pattern already used on other sheets.
$wsZubGel.Range("E9").Value = "((P_Art * Menge_L) + (P_Art_S1 * Menge_S1) + (P_Art_S2 * Menge_S2))* ( 1 - (p_rabatt / 100)"

# --- restore the cell selections as left by the editing session -----------
$wsBrixGelStab = $wb.Sheets.Item("Brix_Gel_Stab")
[void]$wsBrixGelStab.Activate()
[void]$wsBrixGelStab.Range("E12").Select()

[void]$wsZubGel.Activate()
[void]$wsZubGel.Range("D16").Select()
